$wb = $excel.ActiveWorkbook

# --- Create the new "Slovakia" sheet by copying "Portugal" (same layout/styles),
#     placing it immediately after Portugal in the tab order. ---
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Copy([System.Reflection.Missing]::Value, $portugal) | Out-Null
$slovakia = $wb.Worksheets.Item($portugal.Index + 1)
$slovakia.Name = "Slovakia"

# --- Fill in the market-specific values for Slovakia. ---
$slovakia.Range("B2").Value = "Slovakia Market"
$slovakia.Range("B4").Value = "NGC-2930/T3222/T3221/T3225"

# --- Slovakia's printer list does not include "PROFILE Communicator", so remove
#     that row (shifts the remaining printer rows up by one). ---
$slovakia.Rows.Item(14).Delete() | Out-Null

# --- Those rows were taller (ht 28.8) on the Portugal sheet; Slovakia uses the
#     default row height instead. ---
$slovakia.Rows.Item(3).AutoFit() | Out-Null
$slovakia.Rows.Item(4).AutoFit() | Out-Null
$slovakia.Rows.Item(5).AutoFit() | Out-Null

# --- Update Portugal's own selection to a full-sheet selection. ---
$portugal.Cells.Select() | Out-Null

# --- Make Slovakia the active sheet/tab with B4 selected. ---
$slovakia.Activate()
$slovakia.Range("B4").Select() | Out-Null
